$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for specific rows, per repull of data
$ws.Range("F2").Value = 1
$ws.Range("F3").Value = -2
$ws.Range("F4").Value = 3
$ws.Range("F5").Value = -4
$ws.Range("F7").Value = -8
$ws.Range("F10").Value = -2
$ws.Range("F11").Value = 4
$ws.Range("F14").Value = -2
$ws.Range("F23").Value = -2
